$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 14: new test case #13 (ComposedFoodItem.create, isImport = false) ---
$ws.Range("A14").Value = 13
$ws.Range("B14").Value = "CoreDataTests.ComposedFoodItemBehavior"
$ws.Range("C14").Value = "ComposedFoodItem.create(from composedFoodItemVM: ComposedFoodItemViewModel, _ isImport: Bool)"
$ws.Range("D14").Value = "isImport = false"
$ws.Range("E14").Value = "ComposedFoodItemVM with ingredients, which already exist as FoodItems in the DB"
$ws.Range("F14").Value = "ComposedFoodItem, its related FoodItem, all related Ingredients and their related FoodItems are created"
$ws.Range("G14").Value = "yes"
$ws.Rows.Item(14).RowHeight = 34

# --- Row 15: new test case #14 (ComposedFoodItem.update) ---
$ws.Range("A15").Value = 14
$ws.Range("B15").Value = "CoreDataTests.ComposedFoodItemBehavior"
$ws.Range("C15").Value = "ComposedFoodItem.update(_ composedFoodItemVM: ComposedFoodItemViewModel)"
$ws.Range("E15").Value = "A fully available ComposedFoodItem with related FoodItem, Ingredients and their FoodItems"
$ws.Rows.Item(15).RowHeight = 34

# --- Update view: top-left visible cell scrolled to A2, selection moved to F15 ---
try {
    $excel.ActiveWindow.ScrollRow = 2
    $excel.ActiveWindow.ScrollColumn = 1
} catch {}
[void]$ws.Range("F15").Select()
